$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "Chalcone Synthase 1"
$ws.Range("D8").Value = "Chalcone Synthase 1"
$ws.Range("D9").Value = "Chalcone Synthase"
$ws.Range("D17").Value = "Chalcone Synthase J"
$ws.Range("D18").Value = "Chalcone Synthase 2"

$ws.Range("D17").Select()
